{"js": "// Update the date heading and the 25 \"two-digit \u00f7 one-digit\" answer\n// cells in the practice-sheet table to the next day's generated values.\n//\n// The document body is: one centered paragraph with the date, followed\n// by a single 5-column table whose data rows (every 4th row, the other\n// three being blank spacer rows) hold the answer strings. We replace\n// the date paragraph's text, then walk the table's non-empty cells in\n// document order and replace their text according to the recorded\n// old -> new mapping (applied positionally, not via text search, so\n// the fact that one new value happens to equal another cell's old\n// value can't cause a double-replace).\n\nconst oldToNewDate = {\n  \"2025-01-28 Tuesday\": \"2025-01-29 Wednesday\",\n};\n\n// Position-ordered (row-major, left-to-right, top-to-bottom over the\n// non-blank rows) old -> new answers for the 25 table cells.\nconst cellReplacements = [\n  [\"29\u00f77=4, 1\", \"58\u00f78=7, 2\"],\n  [\"90\u00f76=15, 0\", \"20\u00f76=3, 2\"],\n  [\"60\u00f75=12, 0\", \"54\u00f72=27, 0\"],\n  [\"84\u00f77=12, 0\", \"53\u00f74=13, 1\"],\n  [\"84\u00f75=16, 4\", \"14\u00f78=1, 6\"],\n  [\"70\u00f76=11, 4\", \"23\u00f75=4, 3\"],\n  [\"33\u00f72=16, 1\", \"74\u00f73=24, 2\"],\n  [\"36\u00f77=5, 1\", \"60\u00f74=15, 0\"],\n  [\"59\u00f76=9, 5\", \"52\u00f76=8, 4\"],\n  [\"83\u00f72=41, 1\", \"69\u00f75=13, 4\"],\n  [\"62\u00f78=7, 6\", \"62\u00f79=6, 8\"],\n  [\"22\u00f74=5, 2\", \"85\u00f75=17, 0\"],\n  [\"84\u00f78=10, 4\", \"38\u00f77=5, 3\"],\n  [\"76\u00f79=8, 4\", \"12\u00f79=1, 3\"],\n  [\"84\u00f74=21, 0\", \"88\u00f74=22, 0\"],\n  [\"46\u00f77=6, 4\", \"84\u00f77=12, 0\"],\n  [\"76\u00f78=9, 4\", \"63\u00f79=7, 0\"],\n  [\"12\u00f74=3, 0\", \"87\u00f75=17, 2\"],\n  [\"61\u00f73=20, 1\", \"60\u00f77=8, 4\"],\n  [\"64\u00f77=9, 1\", \"60\u00f79=6, 6\"],\n  [\"36\u00f76=6, 0\", \"51\u00f77=7, 2\"],\n  [\"10\u00f75=2, 0\", \"28\u00f77=4, 0\"],\n  [\"98\u00f72=49, 0\", \"19\u00f72=9, 1\"],\n  [\"70\u00f73=23, 1\", \"50\u00f76=8, 2\"],\n  [\"40\u00f74=10, 0\", \"71\u00f74=17, 3\"],\n];\n\nconst body = context.document.body;\n\n// --- 1) The date paragraph (first paragraph of the body, outside the\n// table) -----------------------------------------------------------\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (const para of paragraphs.items) {\n  para.load(\"text\");\n}\nawait context.sync();\n\nfor (const para of paragraphs.items) {\n  const replacement = oldToNewDate[para.text];\n  if (replacement !== undefined) {\n    para.getRange().insertText(replacement, \"Replace\");\n  }\n}\nawait context.sync();\n\n// --- 2) The table's answer cells -----------------------------------\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nfor (const row of rows.items) {\n  row.cells.load(\"items\");\n}\nawait context.sync();\n\n// Collect cells in row-major document order.\nconst allCells = [];\nfor (const row of rows.items) {\n  for (const cell of row.cells.items) {\n    allCells.push(cell);\n  }\n}\n\nfor (const cell of allCells) {\n  cell.body.load(\"text\");\n}\nawait context.sync();\n\nlet dataIdx = 0;\nfor (const cell of allCells) {\n  if (cell.body.text === \"\") {\n    continue; // blank spacer-row cell, skip\n  }\n  const [, newValue] = cellReplacements[dataIdx];\n  dataIdx++;\n  const para = cell.body.paragraphs.getFirst();\n  para.getRange().insertText(newValue, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Update the date heading and the 25 \"two-digit \u00f7 one-digit\" answer\n# cells in the practice-sheet table to the next day's generated values.\n#\n# The document body is: one centered paragraph with the date, followed\n# by a single 5-column table whose data rows (every 4th row, the other\n# three being blank spacer rows) hold the answer strings. We replace\n# the date paragraph's Range.Text, then walk the table cells in\n# document order (row-major) and replace the non-blank ones according\n# to the recorded old -> new mapping, applied positionally so the fact\n# that one new value happens to equal another cell's old value can't\n# cause a double-replace.\n\n$d = $word.ActiveDocument\n\n# --- 1) The date paragraph (first paragraph of the body, outside the\n# table) --------------------------------------------------------------\n$d.Paragraphs(1).Range.Text = \"2025-01-29 Wednesday\"\n\n# --- 2) The table's answer cells --------------------------------------\n# Position-ordered (row-major, left-to-right, top-to-bottom over the\n# non-blank rows) new answers for the 25 table cells.\n$newValues = @(\n    \"58\u00f78=7, 2\",\n    \"20\u00f76=3, 2\",\n    \"54\u00f72=27, 0\",\n    \"53\u00f74=13, 1\",\n    \"14\u00f78=1, 6\",\n    \"23\u00f75=4, 3\",\n    \"74\u00f73=24, 2\",\n    \"60\u00f74=15, 0\",\n    \"52\u00f76=8, 4\",\n    \"69\u00f75=13, 4\",\n    \"62\u00f79=6, 8\",\n    \"85\u00f75=17, 0\",\n    \"38\u00f77=5, 3\",\n    \"12\u00f79=1, 3\",\n    \"88\u00f74=22, 0\",\n    \"84\u00f77=12, 0\",\n    \"63\u00f79=7, 0\",\n    \"87\u00f75=17, 2\",\n    \"60\u00f77=8, 4\",\n    \"60\u00f79=6, 6\",\n    \"51\u00f77=7, 2\",\n    \"28\u00f77=4, 0\",\n    \"19\u00f72=9, 1\",\n    \"50\u00f76=8, 2\",\n    \"71\u00f74=17, 3\"\n)\n\n$t = $d.Tables(1)\n$idx = 0\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    for ($c = 1; $c -le $t.Columns.Count; $c++) {\n        $cell = $t.Cell($r, $c)\n        $cellText = $cell.Range.Text -replace \"[\\r\\a]\", \"\"\n        if ($cellText -eq \"\") {\n            continue\n        }\n        $cell.Range.Text = $newValues[$idx]\n        $idx++\n    }\n}\n"}
